# Weekly update: a new price record (week of 2023-12-05) is inserted as a
# new row right before the existing row 1026, shifting every subsequent
# row down by one (old row 1104 becomes new row 1105; dimension grows
# from A1:R1104 to A1:R1105).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 1026, pushing rows 1026-1104 down to 1027-1105.
$ws.Rows.Item(1026).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Range("A1026").Value = 5
$ws.Range("B1026").Value = "Macroferia Regional de Talca"
$ws.Range("C1026").Value = "Maule"
$ws.Range("D1026").Value2 = 45265
$ws.Range("E1026").Value = 7
$ws.Range("F1026").Value = 100112002
$ws.Range("G1026").Value = "Pimiento"
$ws.Range("H1026").Value = "Cuatro cascos verde"
$ws.Range("I1026").Value = "Primera"
$ws.Range("J1026").Value = 200
$ws.Range("K1026").Value = 15000
$ws.Range("L1026").Value = 15000
$ws.Range("M1026").Value = 15000
$ws.Range("N1026").Value = '$/caja 15 kilos'
$ws.Range("O1026").Value = "Región del Maule"
$ws.Range("P1026").Value = 1000
$ws.Range("Q1026").Value = 15
$ws.Range("R1026").Value = "Hortaliza"
